$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new ratings row (No=37, Name=Vons, Gender=M, Rating=24) right
# below the existing data, mirroring the formatting used by the rows above.
$ws.Range("A37:E37").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "Vons"
$ws.Range("C38").Value = "M"
$ws.Range("D38").Clear()
$ws.Range("E38").Value = 24

$ws.Rows.Item(38).RowHeight = 15.75

# Match the author's last on-screen selection.
$ws.Range("H41").Select() | Out-Null
